$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$style = $ws.Range("D2").Style
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = '70.583.66'
$ws.Range("D2").Style = $style
$ws.Range("E2").Value2 = '  +0.66%  '
$style = $ws.Range("D3").Style
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = '3.622.74'
$ws.Range("D3").Style = $style
$ws.Range("E3").Value2 = '  +2.25%  '
$style = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = '1.00'
$ws.Range("D4").Style = $style
$ws.Range("E4").Value2 = '  +0.16%  '
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = '604.51'
$ws.Range("D5").Style = $style
$ws.Range("E5").Value2 = '  +0.07%  '
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = '196.59'
$ws.Range("D6").Style = $style
$ws.Range("E6").Value2 = '  -0.25%  '
$ws.Range("E7").Value2 = '  +0.08%  '
$style = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value2 = '1.00'
$ws.Range("D8").Style = $style
$ws.Range("E8").Value2 = '  +0.07%  '
$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = '0.207'
$ws.Range("D9").Style = $style
$ws.Range("E9").Value2 = '  -1.23%  '
$ws.Range("E10").Value2 = '  -1.13%  '
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = '53.69'
$ws.Range("D11").Style = $style
$ws.Range("E11").Value2 = '  -0.68%  '
$ws.Range("E12").Value2 = '  +0.03%  '
$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = '9.58'
$ws.Range("D13").Style = $style
$ws.Range("E13").Value2 = '  +0.24%  '
$style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = '4.195.69'
$ws.Range("D14").Style = $style
$ws.Range("E14").Value2 = '  +2.07%  '
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = '12.99'
$ws.Range("D15").Style = $style
$ws.Range("E15").Value2 = '  +2.16%  '
$style = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = '596.18'
$ws.Range("D16").Style = $style
$ws.Range("E16").Value2 = '  -1.34%  '
$style = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = '70.662.22'
$ws.Range("D17").Style = $style
$ws.Range("E17").Value2 = '  +0.65%  '
$style = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = '3.634.23'
$ws.Range("D18").Style = $style
$ws.Range("E18").Value2 = '  +2.44%  '
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value2 = '19.08'
$ws.Range("D19").Style = $style
$ws.Range("E19").Value2 = '  -0.90%  '
$ws.Range("E20").Value2 = '  +1.45%  '
$ws.Range("E21").Value2 = '  +0.13%  '
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = '17.83'
$ws.Range("D22").Style = $style
$ws.Range("E22").Value2 = '  -1.06%  '
$ws.Range("E23").Value2 = '  -2.03%  '
$style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = '101.75'
$ws.Range("D24").Style = $style
$ws.Range("E25").Value2 = '  +0.22%  '
$ws.Range("E26").Value2 = '  -3.83%  '
$ws.Range("E27").Value2 = '  -2.23%  '
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = '9.63'
$ws.Range("D28").Style = $style
$ws.Range("E28").Value2 = '  -0.17%  '
$style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = '33.84'
$ws.Range("D29").Style = $style
$ws.Range("E29").Value2 = '  +0.12%  '
$style = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = '4.67'
$ws.Range("D30").Style = $style
$ws.Range("E30").Value2 = '  +6.29%  '
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = '7.25'
$ws.Range("D31").Style = $style
$ws.Range("E31").Value2 = '  +1.43%  '
$ws.Range("E32").Value2 = '  -2.80%  '
$ws.Range("E33").Value2 = '  +1.82%  '
$style = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = '63.64'
$ws.Range("D34").Style = $style
$ws.Range("E34").Value2 = '  +0.35%  '
$ws.Range("E35").Value2 = '  +6.18%  '
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = '3.908.09'
$ws.Range("D36").Style = $style
$ws.Range("E36").Value2 = '  +3.40%  '
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = '542.34'
$ws.Range("D37").Style = $style
$ws.Range("E37").Value2 = '  +10.63%  '
$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = '3.13'
$ws.Range("D38").Style = $style
$ws.Range("E38").Value2 = '  +1.56%  '
$ws.Range("E39").Value2 = '  +0.00%  '
$style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = '37.05'
$ws.Range("D40").Style = $style
$ws.Range("E40").Value2 = '  +0.65%  '
$ws.Range("E41").Value2 = '  -1.39%  '
$ws.Range("E42").Value2 = '  -4.31%  '
$ws.Range("E43").Value2 = '  -0.28%  '
$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = '0.0458'
$ws.Range("D44").Style = $style
$ws.Range("E44").Value2 = '  -0.08%  '
$style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = '3.41'
$ws.Range("D45").Style = $style
$ws.Range("E45").Value2 = '  +3.73%  '
$ws.Range("E46").Value2 = '  +0.33%  '
$style = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = '0.141'
$ws.Range("D47").Style = $style
$ws.Range("E47").Value2 = '  +0.22%  '
$ws.Range("E48").Value2 = '  -0.85%  '
$ws.Range("E49").Value2 = '  -0.12%  '
$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = '0.000252'
$ws.Range("D50").Style = $style
$ws.Range("E50").Value2 = '  +0.84%  '
$ws.Range("E51").Value2 = '  +1.24%  '
